$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "71×83=5893"
$t.Cell(1,2).Range.Text = "92×95=8740"
$t.Cell(1,3).Range.Text = "78×37=2886"
$t.Cell(1,4).Range.Text = "51×29=1479"
$t.Cell(1,5).Range.Text = "26×98=2548"
$t.Cell(2,1).Range.Text = "30×83=2490"
$t.Cell(2,2).Range.Text = "68×79=5372"
$t.Cell(2,3).Range.Text = "81×23=1863"
$t.Cell(2,4).Range.Text = "85×10=850"
$t.Cell(2,5).Range.Text = "66×85=5610"
$t.Cell(3,1).Range.Text = "52×72=3744"
$t.Cell(3,2).Range.Text = "20×72=1440"
$t.Cell(3,3).Range.Text = "55×49=2695"
$t.Cell(3,4).Range.Text = "99×18=1782"
$t.Cell(3,5).Range.Text = "88×80=7040"
$t.Cell(4,1).Range.Text = "42×39=1638"
$t.Cell(4,2).Range.Text = "27×98=2646"
$t.Cell(4,3).Range.Text = "36×98=3528"
$t.Cell(4,4).Range.Text = "60×55=3300"
$t.Cell(4,5).Range.Text = "15×48=720"
$t.Cell(5,1).Range.Text = "62×85=5270"
$t.Cell(5,2).Range.Text = "52×18=936"
$t.Cell(5,3).Range.Text = "72×24=1728"
$t.Cell(5,4).Range.Text = "20×94=1880"
$t.Cell(5,5).Range.Text = "35×34=1190"
$t.Cell(6,1).Range.Text = "43×62=2666"
$t.Cell(6,2).Range.Text = "61×12=732"
$t.Cell(6,3).Range.Text = "15×31=465"
$t.Cell(6,4).Range.Text = "88×54=4752"
$t.Cell(6,5).Range.Text = "78×50=3900"
$t.Cell(7,1).Range.Text = "86×53=4558"
$t.Cell(7,2).Range.Text = "78×96=7488"
$t.Cell(7,3).Range.Text = "29×69=2001"
$t.Cell(7,4).Range.Text = "57×88=5016"
$t.Cell(7,5).Range.Text = "42×41=1722"
$t.Cell(8,1).Range.Text = "80×99=7920"
$t.Cell(8,2).Range.Text = "78×27=2106"
$t.Cell(8,3).Range.Text = "33×41=1353"
$t.Cell(8,4).Range.Text = "59×38=2242"
$t.Cell(8,5).Range.Text = "25×74=1850"
$t.Cell(9,1).Range.Text = "77×55=4235"
$t.Cell(9,2).Range.Text = "42×94=3948"
$t.Cell(9,3).Range.Text = "62×54=3348"
$t.Cell(9,4).Range.Text = "45×90=4050"
$t.Cell(9,5).Range.Text = "72×57=4104"
$t.Cell(10,1).Range.Text = "52×29=1508"
$t.Cell(10,2).Range.Text = "16×37=592"
$t.Cell(10,3).Range.Text = "43×93=3999"
$t.Cell(10,4).Range.Text = "80×60=4800"
$t.Cell(10,5).Range.Text = "80×92=7360"
$t.Cell(11,1).Range.Text = "48×95=4560"
$t.Cell(11,2).Range.Text = "61×86=5246"
$t.Cell(11,3).Range.Text = "91×18=1638"
$t.Cell(11,4).Range.Text = "25×46=1150"
$t.Cell(11,5).Range.Text = "41×12=492"
$t.Cell(12,1).Range.Text = "69×11=759"
$t.Cell(12,2).Range.Text = "40×21=840"
$t.Cell(12,3).Range.Text = "17×33=561"
$t.Cell(12,4).Range.Text = "100×70=7000"
$t.Cell(12,5).Range.Text = "98×17=1666"
$t.Cell(13,1).Range.Text = "77×75=5775"
$t.Cell(13,2).Range.Text = "42×47=1974"
$t.Cell(13,3).Range.Text = "90×86=7740"
$t.Cell(13,4).Range.Text = "61×74=4514"
$t.Cell(13,5).Range.Text = "56×24=1344"
$t.Cell(14,1).Range.Text = "47×30=1410"
$t.Cell(14,2).Range.Text = "92×69=6348"
$t.Cell(14,3).Range.Text = "87×75=6525"
$t.Cell(14,4).Range.Text = "38×38=1444"
$t.Cell(14,5).Range.Text = "76×29=2204"
$t.Cell(15,1).Range.Text = "25×32=800"
$t.Cell(15,2).Range.Text = "61×58=3538"
$t.Cell(15,3).Range.Text = "16×48=768"
$t.Cell(15,4).Range.Text = "51×34=1734"
$t.Cell(15,5).Range.Text = "68×69=4692"
$t.Cell(16,1).Range.Text = "96×52=4992"
$t.Cell(16,2).Range.Text = "59×33=1947"
$t.Cell(16,3).Range.Text = "48×91=4368"
$t.Cell(16,4).Range.Text = "61×87=5307"
$t.Cell(16,5).Range.Text = "36×57=2052"
$t.Cell(17,1).Range.Text = "14×20=280"
$t.Cell(17,2).Range.Text = "64×73=4672"
$t.Cell(17,3).Range.Text = "59×11=649"
$t.Cell(17,4).Range.Text = "60×19=1140"
$t.Cell(17,5).Range.Text = "59×92=5428"
$t.Cell(18,1).Range.Text = "31×29=899"
$t.Cell(18,2).Range.Text = "29×40=1160"
$t.Cell(18,3).Range.Text = "97×95=9215"
$t.Cell(18,4).Range.Text = "16×12=192"
$t.Cell(18,5).Range.Text = "61×90=5490"
$t.Cell(19,1).Range.Text = "26×17=442"
$t.Cell(19,2).Range.Text = "79×45=3555"
$t.Cell(19,3).Range.Text = "29×14=406"
$t.Cell(19,4).Range.Text = "49×84=4116"
$t.Cell(19,5).Range.Text = "13×87=1131"
$t.Cell(20,1).Range.Text = "31×30=930"
$t.Cell(20,2).Range.Text = "13×28=364"
$t.Cell(20,3).Range.Text = "55×96=5280"
$t.Cell(20,4).Range.Text = "51×26=1326"
$t.Cell(20,5).Range.Text = "29×99=2871"
